$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix capitalisation of an existing answer: "svensk farmor" -> "Svensk farmor"
$fixCell = $ws.Cells.Find("svensk farmor")
if ($fixCell) {
    $fixCell.Value2 = "Svensk farmor"
} else {
    $ws.Range("B17").Value2 = "Svensk farmor"
}

# Add the new question/answer pair as a new table row (appended, then the
# table is re-sorted below so it lands in the correct alphabetical spot).
$tbl = $ws.ListObjects.Item(1)
$newRow = $tbl.ListRows.Add()
$newRowRange = $newRow.Range
$newRowRange.Cells.Item(1, 1).Value2 = "I float in the air, bringing joy your way. Vibrant and cheerful, I rise above, Made from old noses, spreading love. What am I?"
$newRowRange.Cells.Item(1, 2).Value2 = "Rød ballong"

# Re-apply the table's existing sort so the new row moves into its sorted
# position (the data is sorted A-Z by the Question column).
$tbl.Sort.Apply()

# Update the active selection to reflect where the user ended up (right
# below the newly added/sorted table row).
$ws.Range("B22").Select()
